$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(42613.757210648146, 20, 55, 41, 55, 18, 12906, 24956, 2888, 370, 274, 39, 9, "Bag"),
    @(42613.884375000001, 20, 58, 39, 58, 36, 7768, 14830, 1704, 216, 146, 23, 13, "Bag"),
    @(42614.883136574077, 4, 46, 45, 46, 28, 14333, 32461, 3720, 426, 412, 38, 15, "Bag"),
    @(42615.884050925924, 12, 52, 46, 52, 27, 14108, 28475, 3223, 372, 329, 42, 16, "Bag")
)

$startRow = 10
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
    $ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
}
